$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wireless_Charging")

# --- Row 7: fill B7:D7 and F7:I7 (E7 stays a standalone formula, matching the
# original author's edit which apparently skipped E7 when re-entering the row) ---
$ws.Range("B7:D7").Formula = '=((($B$11*2*PI())^2*B5)^-1)*1000000000'
$ws.Range("F7:I7").Formula = '=((($B$11*2*PI())^2*F5)^-1)*1000000000'

# --- Row 8: fill B8:I8 ---
$ws.Range("B8:I8").Formula = '=(($B$12*2*PI())^2*B5-1/B7)^-1'

# --- Row 9: fill B9:I9 ---
$ws.Range("B9:I9").Formula = '=(2*PI()*$B$11*B5/B6)'

# The engine's Range.Formula setter sometimes pulls in a donor cell's number
# format (e.g. the scientific-notation format used by F5:I5) when a formula
# references it. Restore the plain centered style (the same style already
# used throughout column F:I on rows 7-9) by copying formatting only from a
# same-style neighbor, so no new style slots get minted and no visible
# format actually changes.
$ws.Range("G6").Copy() | Out-Null
$ws.Range("F7:I9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New cell C11: resonant frequency calculation ---
$ws.Range("C11").Formula = '=1/(2*PI()*SQRT(0.000000047*0.00001215))'

# --- New row 13 ---
$ws.Range("D13").Formula = '=100/1.29'
$ws.Range("E13").Formula = '=1/(2*PI()*SQRT(0.0000126*0.000000047))'

# --- New row 14 ---
$ws.Range("E14").Formula = '=1/(2*PI()*SQRT(0.0000072*0.0000000033))'

# Update the active selection to C12
$ws.Range("C12").Select()
